# "Generate Report for handback" — refresh the handback-completion
# timestamps on the per-locale status sheets (Correspond Handoff
# Datetime / Correspond Handback DateTime columns, D & G) now that the
# handback for these two files has actually completed.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-01-08 11:47:29"
$zhcn.Range("G2").Value = "2016-01-08 11:48:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-01-08 11:47:42"
$dede.Range("G2").Value = "2016-01-08 11:48:39"
